# Nalco aluminium ingot price sheet update.
#
# The source publishes one new price row at the top each time it is
# refreshed: every existing row (2..150) shifts down by one position and
# a brand new row appears at the bottom that carries forward the data
# (and hyperlink) that used to belong to the last row. Row 1 is the
# header and is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 2. This physically shifts rows 2..150
#    down to 3..151 and automatically extends the used range/dimension
#    to A1:F151 (row 151 already ends up with the same text/number
#    values that used to live in row 150).
$ws.Rows.Item(2).Insert()

# 2) Populate the brand-new row 2 with the latest circular's data.
#    The Date/Circular Date columns are stored as plain text (not real
#    dates) in this sheet, so force a Text number format before
#    assigning them to stop Excel from auto-converting the string into
#    a date value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "03-01-2026"

$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 307.25

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-01-2026"

$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

# Re-apply the standard data-row formatting (General number format and
# the same alignment/border every other row uses) so row 2 looks
# identical to the rest of the table instead of keeping the Text format
# tweak above or whatever formatting Insert() guessed at.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$wb.Application.CutCopyMode = 0

# 3) Row-insert in this environment does not reliably carry each cell's
#    hyperlink object down with its row, so rebuild every hyperlink in
#    column F from scratch: each F cell's link target always equals its
#    own displayed text in this sheet, so that is a safe source of
#    truth for all 150 data rows (2..151), including the newly added
#    row 151 which has no hyperlink object yet.
for ($r = 2; $r -le 151; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Hyperlinks.Count -gt 0) {
        $cell.Hyperlinks.Delete()
    }
}
for ($r = 2; $r -le 151; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
}

# Adding hyperlinks switches those cells to the default blue/underlined
# "Hyperlink" style; restore the plain data-row style used by every
# other cell in the table (column C carries that same plain style).
$ws.Range("C3").Copy()
$ws.Range("F2:F151").PasteSpecial(-4122)
$wb.Application.CutCopyMode = 0
